$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Works?" column (C) for the next three upgrade rows
$ws.Range("C9").Value = "Yes"
$ws.Range("C10").Value = "No"
$ws.Range("C11").Value = "No"

# Update the view: zoom to 90% and move the active selection to C12
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("C12").Select()
